# Generate Report for Handback
#
# The localization-status report previously listed
# "62e63023-0051-4a1b-b965-de3d015188b2.md" as "Ready for handoff" (row 7)
# and "ed807c42-2548-4889-a845-e2dc2186f24d.md" as "In Translation" (row 6)
# on the Overview / zh-cn / de-de sheets.
#
# A handback run was generated for "62e63023...md" which failed the
# handback transform, so it now moves up to row 6 with a
# "Handback transform failed" status, an updated handoff timestamp, and a
# new Error Detail message. The "ed807c42...md" entry moves down to row 7,
# keeping its old "In Translation" values untouched.

$wb = $excel.ActiveWorkbook

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L")

# ---------------------------------------------------------------------
# 1) Swap the contents of row 6 and row 7 on every sheet (generic swap
#    across all columns actually used by the table), then patch in the
#    new status/timestamp/error-detail values for the row that now holds
#    the 62e63023 entry.
# ---------------------------------------------------------------------

$sheetNames = @("Overview","zh-cn","de-de")

foreach ($sname in $sheetNames) {
    $ws = $wb.Worksheets.Item($sname)

    # Read current row 6 / row 7 values for every column first (so the
    # swap doesn't clobber values we still need to read).
    $row6 = @{}
    $row7 = @{}
    foreach ($c in $cols) {
        $row6[$c] = $ws.Range($c + "6").Value()
        $row7[$c] = $ws.Range($c + "7").Value()
    }

    # Write the swap: new row 6 = old row 7, new row 7 = old row 6.
    foreach ($c in $cols) {
        $ws.Range($c + "6").Value = $row7[$c]
        $ws.Range($c + "7").Value = $row6[$c]
    }
}

# ---------------------------------------------------------------------
# 2) Patch the now-row-6 (62e63023...) entry with its handback-failure
#    status, refreshed "latest handoff" timestamp and error detail.
# ---------------------------------------------------------------------

# Overview sheet: columns B (zh-cn) and C (de-de) hold the status, column
# D holds the "Latest Handoff Date" shown on the overview.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B6").Value = "Handback transform failed"
$wsOverview.Range("C6").Value = "Handback transform failed"
$wsOverview.Range("D6").Value = "2016-03-21 16:25:51"

# zh-cn sheet: column C is Status, column E is Latest Handoff Datetime,
# column L is the new Error Detail.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C6").Value = "Handback transform failed"
$wsZh.Range("E6").Value = "2016-03-21 16:25:46"
$wsZh.Range("L6").Value = "The handback type mt is not match with handoff type ht."

# de-de sheet: same columns as zh-cn, different timestamp.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C6").Value = "Handback transform failed"
$wsDe.Range("E6").Value = "2016-03-21 16:25:51"
$wsDe.Range("L6").Value = "The handback type mt is not match with handoff type ht."

# ---------------------------------------------------------------------
# 3) Fix up the hyperlinks so row 6 links to 62e63023... and row 7 links
#    to ed807c42..., matching the swapped cell content. The underlying
#    rId/target URL for each hyperlink stays attached to its row
#    position; only the displayed text changes.
# ---------------------------------------------------------------------

$overviewDisplay = @{}
$overviewDisplay['$A$6'] = "62e63023-0051-4a1b-b965-de3d015188b2.md"
$overviewDisplay['$A$7'] = "ed807c42-2548-4889-a845-e2dc2186f24d.md"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($overviewDisplay.ContainsKey($addr)) {
        $h.TextToDisplay = $overviewDisplay[$addr]
    }
}

$zhDisplay = @{}
$zhDisplay['$A$6'] = "62e63023-0051-4a1b-b965-de3d015188b2.md"
$zhDisplay['$D$6'] = "62e63023-0051-4a1b-b965-de3d015188b2.505be5758d63a37bb16447aa0dce246d9b377efc.zh-cn.xlf"
$zhDisplay['$A$7'] = "ed807c42-2548-4889-a845-e2dc2186f24d.md"
$zhDisplay['$D$7'] = "ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.zh-cn.xlf"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($zhDisplay.ContainsKey($addr)) {
        $h.TextToDisplay = $zhDisplay[$addr]
    }
}

$deDisplay = @{}
$deDisplay['$A$6'] = "62e63023-0051-4a1b-b965-de3d015188b2.md"
$deDisplay['$D$6'] = "62e63023-0051-4a1b-b965-de3d015188b2.505be5758d63a37bb16447aa0dce246d9b377efc.de-de.xlf"
$deDisplay['$A$7'] = "ed807c42-2548-4889-a845-e2dc2186f24d.md"
$deDisplay['$D$7'] = "ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.de-de.xlf"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($deDisplay.ContainsKey($addr)) {
        $h.TextToDisplay = $deDisplay[$addr]
    }
}
